$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "AMSIN": append two new interview-history rows (21 and 22)
# ---------------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

# Seed row 21 by duplicating the last existing row (20) so the new row
# inherits the same look & feel (general number format / date-time format)
# as the rest of the table, then overwrite the cell values.
$wsAmsin.Range("A20:G20").Copy()
$wsAmsin.Range("A21:G21").PasteSpecial(-4104)
$wsAmsin.Range("A22:G22").PasteSpecial(-4104)
$excel.CutCopyMode = 0

$wsAmsin.Range("A21").Value = "'2021-10-26"
$wsAmsin.Range("B21").Value = 44495.64975144676
$wsAmsin.Range("C21").Value = "152_fstcycle"
$wsAmsin.Range("D21").Value = 155
$wsAmsin.Range("E21").Value = 150
$wsAmsin.Range("F21").Value = 5
$wsAmsin.Range("G21").Value = 4.84

$wsAmsin.Range("A22").Value = "'2021-10-28"
$wsAmsin.Range("B22").Value = 44497.3907918287
$wsAmsin.Range("C22").Value = "152_fnlrgrsn"
$wsAmsin.Range("D22").Value = 155
$wsAmsin.Range("E22").Value = 154
$wsAmsin.Range("F22").Value = 1
$wsAmsin.Range("G22").Value = 3.8

# Restore the general formatting on column A / B after the value write
# (typing a value can re-derive a number format) and the date-time format
# used throughout the report on column B.
$wsAmsin.Range("A20").Copy()
$wsAmsin.Range("A21").PasteSpecial(-4122)
$wsAmsin.Range("A22").PasteSpecial(-4122)
$wsAmsin.Range("B20").Copy()
$wsAmsin.Range("B21").PasteSpecial(-4122)
$wsAmsin.Range("B22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Sheet "AMS": row 19 picks up the table's shared formatting (it previously
# had none), and a new row 20 is appended with the next interview run.
# ---------------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

$wsAms.Range("A18:G18").Copy()
$wsAms.Range("A19:G19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$wsAms.Range("A20").Value = "'2021-10-28"
$wsAms.Range("B20").Value = 44497.87484470283
$wsAms.Range("C20").Value = "152_livetest"
$wsAms.Range("D20").Value = 155
$wsAms.Range("E20").Value = 153
$wsAms.Range("F20").Value = 2
$wsAms.Range("G20").Value = 4.16

$wsAms.Range("B19").Copy()
$wsAms.Range("B20").PasteSpecial(-4122)
$excel.CutCopyMode = 0
